# Add a new worksheet "ODI Batting Extra" (sheetId=3, rId3) after "ODI Batting",
# containing the MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns with one data row for match 4690.

$wb = $excel.ActiveWorkbook

# Keep a handle to the original first sheet so we can re-activate it at the end
# (adding/activating the new sheet shouldn't change which tab is active).
$firstSheet = $wb.Worksheets.Item(1)

# Insert the new sheet after the last existing sheet ("ODI Batting") so it
# becomes the 3rd / final tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Match the bold / centered / bordered header styling used on the other sheets
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row - MATCH_CODE must stay textual ("4690"), not numeric.
# BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are left blank.
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "4690"
$newSheet.Range("A2").Style = "Normal"
$newSheet.Range("F2").Value = "NO"

# Restore original active sheet/selection
$firstSheet.Activate()

Write-Host "Added sheet 'ODI Batting Extra' with match 4690 summary row"
